# Regenerate the "K" column (column G) of the save_data sheet.
# This corresponds to switching the source metric used to populate
# column G from the old "Strike#" based figure to the newly computed
# "K" value (std/mean derived), and writing the recalculated s_vals
# back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number (sheet row) -> new value for column G ("K")
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    6  = 1
    7  = 0
    8  = 2
    9  = 2
    10 = 2
    11 = 3
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 5
    20 = 5
    21 = 1
    22 = 2
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 2
    30 = 1
    31 = 0
    32 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 1
    39 = 3
    40 = 2
    41 = 0
    42 = 0
    43 = 1
    44 = 3
    45 = 0
    46 = 3
    47 = 2
    48 = 2
    49 = 0
    50 = 1
    52 = 1
    53 = 6
    54 = 2
    55 = 2
    56 = 2
    57 = 2
    59 = 3
    60 = 0
    61 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
